$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.865.24"
$ws.Range("E2").Value = "  +5.69%  "

$ws.Range("D3").Value = "2.248.78"
$ws.Range("E3").Value = "  +4.13%  "

$ws.Range("E4").Value = "  -0.06%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "229.85"
$ws.Range("E5").Value = "  +1.09%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.627"
$ws.Range("E6").Value = "  -0.21%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "61.13"
$ws.Range("E7").Value = "  -3.53%  "

$ws.Range("E8").Value = "  +0.05%  "

$ws.Range("E9").Value = "  +3.68%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "58.66"
$ws.Range("E10").Value = "  +1.13%  "

$ws.Range("E11").Value = "  +4.64%  "

$ws.Range("E12").Value = "  +0.30%  "

$ws.Range("D13").Value = "2.583.50"
$ws.Range("E13").Value = "  +4.01%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "15.80"
$ws.Range("E14").Value = "  -0.73%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "21.83"
$ws.Range("E15").Value = "  -0.39%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.804"
$ws.Range("E16").Value = "  -0.18%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.62"
$ws.Range("E17").Value = "  +2.69%  "

$ws.Range("D18").Value = "2.246.54"
$ws.Range("E18").Value = "  +3.99%  "

$ws.Range("D19").Value = "41.831.92"
$ws.Range("E19").Value = "  +5.79%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "73.23"
$ws.Range("E20").Value = "  +2.08%  "

$ws.Range("B21").Value = "ShibaInu"
$ws.Range("C21").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D21").Value = "0.0₃0897"
$ws.Range("E21").Value = "  +6.49%  "

$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.12"
$ws.Range("E22").Value = "  +0.33%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "249.54"
$ws.Range("E23").Value = "  +9.78%  "

$ws.Range("E24").Value = "  +0.04%  "

$ws.Range("E25").Value = "  -0.33%  "

$ws.Range("E26").Value = "  -2.26%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.67"
$ws.Range("E27").Value = "  +2.46%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "168.33"
$ws.Range("E28").Value = "  -2.53%  "

$ws.Range("E29").Value = "  +3.09%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "20.15"
$ws.Range("E30").Value = "  +2.31%  "

$ws.Range("E31").Value = "  +1.90%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.78"
$ws.Range("E32").Value = "  +3.74%  "

$ws.Range("E33").Value = "  +0.73%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.23"
$ws.Range("E34").Value = "  +12.09%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.75"
$ws.Range("E35").Value = "  +3.85%  "

$ws.Range("E36").Value = "  +1.77%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.81"
$ws.Range("E37").Value = "  +5.83%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.70"
$ws.Range("E38").Value = "  -3.49%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.39"
$ws.Range("E39").Value = "  -0.16%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.000250"
$ws.Range("E40").Value = "  +40.46%  "

$ws.Range("B41").Value = "FTXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.05"
$ws.Range("E41").Value = "  -3.73%  "

$ws.Range("B42").Value = "BinanceUSD"
$ws.Range("C42").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.00"
$ws.Range("E42").Value = "  +0.06%  "

$ws.Range("E43").Value = "  +4.72%  "

$ws.Range("E44").Value = "  +12.89%  "

$ws.Range("B45").Value = "Cronos"
$ws.Range("C45").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0983"
$ws.Range("E45").Value = "  +6.92%  "

$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "100.30"
$ws.Range("E46").Value = "  -1.47%  "

$ws.Range("D47").Value = "1.484.68"
$ws.Range("E47").Value = "  -1.75%  "

$ws.Range("E48").Value = "  -3.76%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "16.51"
$ws.Range("E49").Value = "  -5.21%  "

$ws.Range("E50").Value = "  +0.63%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.79"
$ws.Range("E51").Value = "  -0.56%  "

